$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.165.74"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.907.67"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.70"
$ws.Range("E5").Value = "  +0.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9994"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5241"
$ws.Range("E7").Value = "  +2.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3781"
$ws.Range("E8").Value = "  +3.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07304"
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.27"
$ws.Range("E10").Value = "  +2.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8994"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07676"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.887.75"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.81"
$ws.Range("E14").Value = "  -0.15%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.251"
$ws.Range("E15").Value = "  +0.23%  "
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008545"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.54"
$ws.Range("E18").Value = "  +2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.230.04"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.090"
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.133.31"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.63"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.443"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.319"
$ws.Range("E25").Value = "  +10.99%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "145.85"
$ws.Range("E26").Value = "  -1.73%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("E28").Value = "  -3.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.85"
$ws.Range("E29").Value = "  +1.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.960"
$ws.Range("E30").Value = "  +4.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.810"
$ws.Range("E31").Value = "  +1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09208"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05075"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.247"
$ws.Range("E34").Value = "  +7.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7812"
$ws.Range("E35").Value = "  +4.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.992"
$ws.Range("E36").Value = "  +0.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.306"
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.607"
$ws.Range("E38").Value = "  +3.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5685"
$ws.Range("E39").Value = "  +0.77%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01994"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.029"
$ws.Range("E42").Value = "  +5.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.632"
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "118.62"
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1523"
$ws.Range("E45").Value = "  +2.94%  "
$ws.Range("E46").Value = "  +2.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.25"
$ws.Range("E47").Value = "  +0.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9995"
$ws.Range("E48").Value = "  -0.05%  "
$ws.Range("E49").Value = "  +2.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.42"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("E51").Value = "  +1.72%  "
